$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 1.533500966413222
$ws.Range("R2").Value = 13.801508697719
$ws.Range("S2").Value = 0.0662600404061536
$ws.Range("T2").Value = 0.06626004040615362

# Row 3 updates
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("R3").Value = 88.21146330053101
$ws.Range("S3").Value = 0.4234968256437875
$ws.Range("T3").Value = 0.4234968256437876

# Row 4 updates
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("S4").Value = 0.5102431339500588
$ws.Range("T4").Value = 0.5102431339500588
